$wb = $excel.ActiveWorkbook

# --- workbook.xml: ink annotation flag (best-effort; matches the author's
# toggle of Review > Ink > Show annotations off) ---
$wb.DisplayInkComments = $false

$ws4 = $wb.Worksheets.Item("EPFImporter")

# Make room for the new "ConnectionPool" / "EPFIngesterQueue" design notes
# block: rows 17-30 become a blank gap (14 rows), shoving the old
# "EPFIngester" sub-tree (parser / db-util notes) down by 14 rows.
$ws4.Rows("17:30").Insert()

# Fill in the new content in the same order the notes were originally
# typed (keeps the shared-string table append order lined up with the
# authored workbook).
$ws4.Range("A16").Value = "ConnectionPool"
$ws4.Range("A23").Value = "EPFIngesterQueue"
$ws4.Range("B26").Value = "add(fileName)"
$ws4.Range("B24").Value = "setQueueSize()"
$ws4.Range("I25").Value = "ExecutorService executor = Executors.newFixedThreadPool(NTHREDS);"
$ws4.Range("I23").Value = "This is originally intended to use the Executors.newFixedTheadPool() logic"
$ws4.Range("I24").Value = "This could also be implemented using Quartz scheduler threads"
$ws4.Range("C27").Value = "The fileName is all an EPFIngester needs to process"
$ws4.Range("C28").Value = "And a connection Pool"
$ws4.Range("B17").Value = "Use Apache's org.apache.commons.dbcp"
$ws4.Range("B18").Value = "This allows the designation of the minimum and maximum number of concurrent connections."
$ws4.Range("B19").Value = "It also allows for virtual opening and closing of connections where ""closing"" merely returns the connection to the pool."
$ws4.Range("B20").Value = "This also means that the underlying connector needs to open and close connections whenever it performs a single SQL statement or a small block of statements"
$ws4.Range("B25").Value = "setConnectionPool()"
$ws4.Range("B49").Value = "isTableExists()"
$ws4.Range("B48").Value = "getTableColumnCount()"

# Restore the "EPFIngester" header (now at row 30) that used to live at A16
# (re-uses the existing shared string, no new entry added).
$ws4.Range("A30").Value = "EPFIngester"

# Update EPFImporter sheet view (scrolled/selected further down to show the
# newly-added rows)
$win = $wb.Windows.Item(1)
$ws4.Select() | Out-Null
$ws4.Range("B49").Select() | Out-Null
$win.ScrollRow = 29
$win.ScrollColumn = 1

# --- Sheet3: selection moved from A94 to D102 ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("D102").Select() | Out-Null

# Restore EPFImporter as the active/visible tab (it was the active sheet
# before this edit).
$ws4.Select() | Out-Null
$ws4.Range("B49").Select() | Out-Null
